$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: quantity (F) and total value (H) adjustments
$ws.Range("F27").Value = 1199.000
$ws.Range("H27").Value = 12327.85

# Row 43
$ws.Range("F43").Value = 1853.000
$ws.Range("H43").Value = 48266.39

# Row 47
$ws.Range("F47").Value = 465.000
$ws.Range("H47").Value = 64051.42

# Row 53
$ws.Range("F53").Value = 3217.000
$ws.Range("H53").Value = 12487.34

# Row 108
$ws.Range("F108").Value = 10154.800
$ws.Range("H108").Value = 23050.60

# Row 118
$ws.Range("F118").Value = 2135.000
$ws.Range("H118").Value = 5359.62
